$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch workbook calculation from manual to automatic (removes calcMode="manual").
$excel.Calculation = -4105

# Insert a new row 17 (old row 17 "waitForPageToRender" shifts down to row 18),
# then fill it with the new "textBoxShouldHaveValue" entry.
$ws.Rows(17).Insert()

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Text Field"
$ws.Range("C17").Value = "textBoxShouldHaveValue"
$ws.Range("D17").Value = "Accepts two parameters @locator and @testData. It gets the text from textBox and validates against the @testData provided. If the validation fails testing should still continue"
$ws.Range("C17:D17").WrapText = $true

# The shifted-down row (now row 18) keeps its data but its serial number bumps 16 -> 17.
$ws.Range("A18").Value = 17

# Append a brand-new row 19 with the "elementShouldNotBePresent" entry.
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Any"
$ws.Range("C19").Value = "elementShouldNotBePresent"
$ws.Range("D19").Value = "Accepts no parameters and verifies element is not available in DOM. Returns true if element is not available in DOM"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("D19").WrapText = $true

# Move the view/selection down to the newly-added last row, matching the new extent.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("D19").Select()
